# Update the "Source File" column (E) values in Sheet1 to reflect the
# renamed/regenerated reference data file.
#
# Old value: Y4_B2526_General_&_Special_Surgery_1_B1_reference_data_D26112025T134028.xlsx
# New value: Y4_B2526_General_&_Special_surgery_1_B1_reference_data_D23122025T104608.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldValue = "Y4_B2526_General_&_Special_Surgery_1_B1_reference_data_D26112025T134028.xlsx"
$newValue = "Y4_B2526_General_&_Special_surgery_1_B1_reference_data_D23122025T104608.xlsx"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 322
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
